$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Four row-pairs had their match data (columns F..V) swapped between the
#    two rows while keeping the row's own Indice (A) and data_partida (E).
# ---------------------------------------------------------------------------
function Swap-MatchData($r1, $r2) {
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $cols) {
        $vals1[$col] = $ws.Range($col + $r1).Value2
        $vals2[$col] = $ws.Range($col + $r2).Value2
    }
    foreach ($col in $cols) {
        $ws.Range($col + $r1).Value2 = $vals2[$col]
        $ws.Range($col + $r2).Value2 = $vals1[$col]
    }
}

Swap-MatchData 58 59
Swap-MatchData 63 65
Swap-MatchData 66 67
Swap-MatchData 68 69

# ---------------------------------------------------------------------------
# 2) Three brand-new match rows were appended at the bottom (181-183),
#    extending the used range from A1:V180 to A1:V183.
#    Copy the formatting of the last existing row (180) down onto the new
#    rows first, so the new cells pick up the same per-column styles
#    (bold/centered index column, date-formatted data_partida column, ...).
# ---------------------------------------------------------------------------
$ws.Range("A180:V180").Copy()
$ws.Range("A181:V183").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{row=181; A=180; E=45283.67708333334; F="Atl. Madrid";   G=1; H="Sevilla";        I=0; J=1.58; K="20/08/2023 10:02"; L=1.59; M="23/12/2023 16:14"; N=4.25; O="20/08/2023 10:02"; P=4.27; Q="23/12/2023 16:14"; R=5.98; S="20/08/2023 10:02"; T=5.99; U="23/12/2023 16:14"; V="https://www.betexplorer.com/football/spain/laliga/atl-madrid-sevilla/IeuL2xY7/"},
    @{row=182; A=181; E=45293.70833333334; F="Getafe";        G=0; H="Rayo Vallecano"; I=2; J=2.22; K="17/12/2024 18:03"; L=2.4;  M="02/01/2024 16:40"; N=3.03; O="17/12/2024 18:03"; P=2.95; Q="02/01/2024 16:58"; R=3.52; S="17/12/2024 18:03"; T=3.57; U="02/01/2024 16:58"; V="https://www.betexplorer.com/football/spain/laliga/getafe-rayo-vallecano/lYM4iq4I/"},
    @{row=183; A=182; E=45293.80208333334; F="Real Sociedad"; G=1; H="Alaves";         I=1; J=1.55; K="17/12/2024 18:03"; L=1.53; M="02/01/2024 19:10"; N=3.92; O="17/12/2024 18:03"; P=3.87; Q="02/01/2024 19:13"; R=5.92; S="17/12/2024 18:03"; T=8.34; U="02/01/2024 19:13"; V="https://www.betexplorer.com/football/spain/laliga/real-sociedad-alaves/GrNqhGfT/"}
)

foreach ($rd in $newRows) {
    $r = $rd.row
    $ws.Range("A$r").Value2 = $rd.A
    $ws.Range("B$r").Value2 = "spain"
    $ws.Range("C$r").Value2 = "laliga"
    $ws.Range("D$r").Value2 = "2023-2024"
    $ws.Range("E$r").Value2 = $rd.E
    $ws.Range("F$r").Value2 = $rd.F
    $ws.Range("G$r").Value2 = $rd.G
    $ws.Range("H$r").Value2 = $rd.H
    $ws.Range("I$r").Value2 = $rd.I
    $ws.Range("J$r").Value2 = $rd.J
    $ws.Range("K$r").Value2 = $rd.K
    $ws.Range("L$r").Value2 = $rd.L
    $ws.Range("M$r").Value2 = $rd.M
    $ws.Range("N$r").Value2 = $rd.N
    $ws.Range("O$r").Value2 = $rd.O
    $ws.Range("P$r").Value2 = $rd.P
    $ws.Range("Q$r").Value2 = $rd.Q
    $ws.Range("R$r").Value2 = $rd.R
    $ws.Range("S$r").Value2 = $rd.S
    $ws.Range("T$r").Value2 = $rd.T
    $ws.Range("U$r").Value2 = $rd.U
    $ws.Range("V$r").Value2 = $rd.V
}
